# Apply the commit's changes to the workbook:
#  - Metadata sheet: Version, Date, Description values updated
#  - Elements sheet: Extension.value[x] Type(s) text updated + column K widened

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.0.0-ballot -> 2.0.0
$meta.Range("B3").Value = "2.0.0"

# Date: 2025-10-15T15:04:22+00:00 -> 2025-10-16T14:47:34+00:00
$meta.Range("B8").Value = "2025-10-16T14:47:34+00:00"

# Description: updated wording referencing PDSm_SimplifiedPublish
$meta.Range("B12").Value = "Pièces jointes liées à l’événement et à l'évaluation. L'extension référence le profil PDSm_SimplifiedPublish."

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.value[x] Type(s) cell (K6): Attachment -> Reference(...)
$elements.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-document-reference)`n"

# Widen column K (Type(s)) to fit the new, much longer value.
# The runtime quantises ColumnWidth to sixths of a character at save time,
# so 78.8 is the input that round-trips closest to the target 79.70703125.
$elements.Columns.Item(11).ColumnWidth = 78.8
